$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add the title block (name / course / professor / "Homework 2")
#    at the very top of the document.
#
#    Each new paragraph is created at the end of the story (so it
#    picks up no inherited paragraph style) and then moved to the
#    front with Cut/Paste, which keeps it free of any inherited
#    pPr/style from neighbouring paragraphs.
# ------------------------------------------------------------------
function InsertAtFront($text, $style) {
    $p = $d.Paragraphs.Add()
    $p.Range.Text = $text
    if ($style) {
        $p.Range.Style = $style
    }
    $p.Range.Cut()
    $front = $d.Range(0, 0)
    $front.Paste()
}

InsertAtFront "Homework 2" "Title"
InsertAtFront "Junggab Son" $null
InsertAtFront "CS422 Machine Learning" $null
InsertAtFront "Brandon Timok" $null

# ------------------------------------------------------------------
# 2. Rewrite the numbered-steps paragraph in the Explanation section.
# ------------------------------------------------------------------
$v = [char]11   # Word's internal "manual line break" (<w:br/>) character

$lines = @(
    "This program implements K-nearest neighbors (KNN) algorithm from scratch using Python. ",
    "It takes two datasets, MNIST_training.csv and MNIST_test.csv, and follows the steps below:",
    "                  ",
    "1. Load the training and test data using pandas.",
    "2. Calculate the Euclidean distance between test and training data using numpy.",
    "3. Finds the K-nearest neighbors and decide the majority class using numpy and Counter.",
    "4. Compares the prediction with the ground truth in the test data using numpy.",
    "5. Computes accuracy by counting correctly and incorrectly classified samples using numpy.",
    "6. Stores the results in a DataFrame and prints it using pandas.",
    "7. Saves the results in a Word document using the python-docx library.",
    "                  "
)

$newText = $v + ($lines -join $v) + $v + $v

$explanationPara = $d.Paragraphs(7)
$explanationPara.Range.Text = $newText

# ------------------------------------------------------------------
# 3. Reformat the Accuracy column of the results table as percentages.
# ------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.Cell(2, 4).Range.Text = "84.0%"
$tbl.Cell(3, 4).Range.Text = "88.0%"
$tbl.Cell(4, 4).Range.Text = "86.0%"
$tbl.Cell(5, 4).Range.Text = "90.0%"
$tbl.Cell(6, 4).Range.Text = "90.0%"

Write-Output "edit complete"
